# Add the new "PV Allocations" worksheet after "Proposed PVs" and populate
# it with the Node / Persistent-Volume mapping table, then restore the
# view/selection state (active sheet, per-sheet selections) that Excel
# would have saved.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new worksheet as the 4th tab (right after "Proposed PVs")
# ---------------------------------------------------------------------
$proposed = $wb.Sheets("Proposed PVs")
$ws = $wb.Worksheets.Add($null, $proposed)
$ws.Name = "PV Allocations"

# ---------------------------------------------------------------------
# 2. Fill in the table data
#    Columns: A = Node, B = backing ES host, C = Persistent Volume
# ---------------------------------------------------------------------
$data = @(
  @("Node",     "Node", "Persistent Volume"),
  @("master-0", "es-2", "es2-vol6"),
  @("master-1", "es-1", "es1-vol5"),
  @("master-2", "es-3", "es3-vol6"),
  @("data-0",   "es-3", "es3-vol7"),
  @("data-1",   "es-7", "es7-vol6"),
  @("data-2",   "es-5", "es5-vol7"),
  @("data-3",   "es-2", "es2-vol7"),
  @("data-4",   "es-6", "es6-vol7"),
  @("data-5",   "es-1", "es1-vol7"),
  @("data-6",   "es-7", "es7-vol5"),
  @("data-7",   "es-2", "es2-vol5"),
  @("data-8",   "es-3", "es3-vol5"),
  @("data-9",   "es-5", "es5-vol6"),
  @("data-10",  "es-1", "es1-vol6"),
  @("data-11",  "es-6", "es6-vol5")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

$lastRow = $data.Length

# ---------------------------------------------------------------------
# 3. Formatting: bold/large header row, regular-size data rows
# ---------------------------------------------------------------------
$ws.Range("A1:C1").Font.Bold = $true
$ws.Range("A1:C1").Font.Size = 24
$ws.Rows(1).RowHeight = 31

$ws.Range("A2:C" + $lastRow).Font.Size = 18

# Column widths (best-fit-like, matching the rest of the workbook's style)
$ws.Columns.Item(1).ColumnWidth = 12.5 - 0.8333333333333333
$ws.Columns.Item(2).ColumnWidth = 10.5 - 0.8333333333333333
$ws.Columns.Item(3).ColumnWidth = 32.83203125 - 0.8333333333333333

# Selection left on the new sheet by the author
$ws.Range("B10").Select()

# ---------------------------------------------------------------------
# 4. Restore the other sheets' selection / active-tab state
# ---------------------------------------------------------------------
$existing = $wb.Sheets("Existing PVs")
$existing.Activate()
$existing.Range("J2").Select()

# "Proposed PVs" ends up the active/selected tab
$proposed.Activate()
